$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the current (pre-edit) full row contents for rows 1195 and 1196,
# since these values will be copied down two rows (to 1197 and 1198) once
# the new rows are inserted, before row 1195/1196 are overwritten with the
# new week's data.
$row1195 = $ws.Range("A1195:R1195").Value2
$row1196 = $ws.Range("A1196:R1196").Value2

# Insert two new blank rows starting at row 1197; this pushes the old
# rows 1197..1285 down to 1199..1287.
$ws.Range("A1197:R1198").EntireRow.Insert()

# The old content that used to live in rows 1195 and 1196 now needs to be
# placed into the newly inserted rows 1197 and 1198 (a straight copy).
$ws.Range("A1197:R1197").Value2 = $row1195
$ws.Range("A1198:R1198").Value2 = $row1196

# Rows 1195 and 1196 now hold the new week of data: only the date (D),
# volume (J), and for row 1196 also the min/max/avg price (K/L/M) and
# Precio $/Kg (P) change.
$ws.Range("D1195").Value2 = 45265
$ws.Range("J1195").Value2 = 1200

$ws.Range("D1196").Value2 = 45265
$ws.Range("J1196").Value2 = 1400
$ws.Range("K1196").Value2 = 6000
$ws.Range("L1196").Value2 = 7000
$ws.Range("M1196").Value2 = 6571
$ws.Range("P1196").Value2 = 1095
